# "Add files via upload" — fix a typo in the username in B2
# ("lautarino" -> "lautarin") and leave the selection on that cell,
# matching the state the sheet was re-saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("B2").Value = "lautarin"

$ws.Range("B2").Select()
